# Insert new paragraphs describing the "Caveat - Weird interaction" section,
# right after the empty paragraph that follows "Chris' explanation..." and
# right before the document's final (trailing) empty paragraph.

$d = $word.ActiveDocument

$newParagraphsXml = @(
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="both"/></w:pPr></w:p>',
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="both"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Caveat – Weird interaction</w:t></w:r></w:p>',
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="both"/></w:pPr></w:p>',
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t>Plotting bars with activations ends up showing that</w:t></w:r></w:p>',
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t>– rACC is only activated in the R condition</w:t></w:r></w:p>',
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t>– amygdalas are activated in both R and C</w:t></w:r></w:p>',
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t>– fusiform is activated in both R and C</w:t></w:r></w:p>',
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="both"/></w:pPr></w:p>',
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t>which is weird. my intuition is that this effect in the t-statistics is driven by the number of trials</w:t></w:r></w:p>',
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:tab/><w:t>length(R) = 2873</w:t></w:r></w:p>',
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:tab/><w:t>length(L) = 797</w:t></w:r></w:p>',
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:tab/><w:t>length(I) = 741</w:t></w:r></w:p>',
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:tab/><w:t>length(C) = 1240</w:t></w:r></w:p>',
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t>but chris says that can''t be the reason.</w:t></w:r></w:p>'
)

# Anchor: the empty paragraph immediately after "Chris' explanation ..."
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptext = $d.Paragraphs.Item($i).Range.Text
    if ($ptext -like "Chris' explanation is that people move the eyes more in that condition/station.*") {
        $anchorIndex = $i + 1
        break
    }
}

if ($anchorIndex -eq -1) {
    throw "Could not locate anchor paragraph ('Chris' explanation...' follower)"
}

foreach ($xml in $newParagraphsXml) {
    $anchor = $d.Paragraphs.Item($anchorIndex)
    $anchorRange = $anchor.Range
    $anchorRange.InsertParagraphAfter() | Out-Null
    $anchorIndex = $anchorIndex + 1
    $newPara = $d.Paragraphs.Item($anchorIndex)
    $newPara.Range.InsertXML($xml) | Out-Null
}

Write-Output "done; paragraphs now: $($d.Paragraphs.Count)"
